# Add a new "From RCSB" column (F) to the worksheet, populating it with
# the same y/n values as the existing "Is model" column (E), and move the
# active selection to F9 as in the edited workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("F1").Value = "From RCSB"

# Data rows - mirror the "Is model" (column E) values for each row
$ws.Range("F2").Value = "y"
$ws.Range("F3").Value = "y"
$ws.Range("F4").Value = "n"
$ws.Range("F5").Value = "n"

# Update the active cell/selection to match the edited workbook
$ws.Range("F9").Select()
